$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45180 -> 45181, i.e. 2023-09-11 -> 2023-09-12) for every data row (2-111).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    if ($current -eq 45180) {
        $cell.Value = 45181
    }
}
